$wb = $excel.ActiveWorkbook

# --- Update localization status text on every sheet that shows it ---
# "Overview" sheet keeps per-language status in columns E (zh-cn) and F (de-de)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F4").Value = "In Translation"

# "zh-cn" and "de-de" detail sheets keep status in column C
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2:C4").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2:C4").Value = "In Translation"

# --- Re-generate (shrink) the status column widths to fit the new, shorter text ---
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
